$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 so a new record ("FY_4.png") can be added,
# pushing the existing TO_*/TD_* rows down by one.
$ws.Rows.Item(36).Insert()

# Re-run of the recognition pipeline changed the computed distance/score
# columns (B-E) for every row, and the recognised-name / verdict columns
# (F-G) for several rows. Rewrite the full data range A1:G45 to match.

# Row 1: AK_1.png
$ws.Cells.Item(1, 1).Value = "AK_1.png"
$ws.Cells.Item(1, 2).Value = 0.542
$ws.Cells.Item(1, 3).Value = 0.002
$ws.Cells.Item(1, 4).Value = 0.971
$ws.Cells.Item(1, 5).Value = 0.267
$ws.Cells.Item(1, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(1, 7).Value = "Salah"

# Row 2: AK_2.png
$ws.Cells.Item(2, 1).Value = "AK_2.png"
$ws.Cells.Item(2, 2).Value = 0.883
$ws.Cells.Item(2, 3).Value = 0.003
$ws.Cells.Item(2, 4).Value = 0.982
$ws.Cells.Item(2, 5).Value = 0.3
$ws.Cells.Item(2, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(2, 7).Value = "Salah"

# Row 3: AK_3.png
$ws.Cells.Item(3, 1).Value = "AK_3.png"
$ws.Cells.Item(3, 2).Value = 0.862
$ws.Cells.Item(3, 3).Value = 0.003
$ws.Cells.Item(3, 4).Value = 0.973
$ws.Cells.Item(3, 5).Value = 0.533
$ws.Cells.Item(3, 6).Value = "Akhlak Kamiswara"
$ws.Cells.Item(3, 7).Value = "Benar"

# Row 4: AK_4.png
$ws.Cells.Item(4, 1).Value = "AK_4.png"
$ws.Cells.Item(4, 2).Value = 0.547
$ws.Cells.Item(4, 3).Value = 0.002
$ws.Cells.Item(4, 4).Value = 0.962
$ws.Cells.Item(4, 5).Value = 0.4
$ws.Cells.Item(4, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(4, 7).Value = "Salah"

# Row 5: AK_5.png
$ws.Cells.Item(5, 1).Value = "AK_5.png"
$ws.Cells.Item(5, 2).Value = 0.5639999999999999
$ws.Cells.Item(5, 3).Value = 0.002
$ws.Cells.Item(5, 4).Value = 0.984
$ws.Cells.Item(5, 5).Value = 0.267
$ws.Cells.Item(5, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(5, 7).Value = "Salah"

# Row 6: MIB_1.png
$ws.Cells.Item(6, 1).Value = "MIB_1.png"
$ws.Cells.Item(6, 2).Value = 1.256
$ws.Cells.Item(6, 3).Value = 0.004
$ws.Cells.Item(6, 4).Value = 0.953
$ws.Cells.Item(6, 5).Value = 0.5669999999999999
$ws.Cells.Item(6, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(6, 7).Value = "Benar"

# Row 7: MIB_2.png
$ws.Cells.Item(7, 1).Value = "MIB_2.png"
$ws.Cells.Item(7, 2).Value = 1.105
$ws.Cells.Item(7, 3).Value = 0.004
$ws.Cells.Item(7, 4).Value = 0.964
$ws.Cells.Item(7, 5).Value = 0.5
$ws.Cells.Item(7, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(7, 7).Value = "Benar"

# Row 8: MIB_3.png
$ws.Cells.Item(8, 1).Value = "MIB_3.png"
$ws.Cells.Item(8, 2).Value = 1.328
$ws.Cells.Item(8, 3).Value = 0.004
$ws.Cells.Item(8, 4).Value = 0.993
$ws.Cells.Item(8, 5).Value = 0.967
$ws.Cells.Item(8, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(8, 7).Value = "Benar"

# Row 9: MIB_4.png
$ws.Cells.Item(9, 1).Value = "MIB_4.png"
$ws.Cells.Item(9, 2).Value = 1.112
$ws.Cells.Item(9, 3).Value = 0.004
$ws.Cells.Item(9, 4).Value = 0.979
$ws.Cells.Item(9, 5).Value = 0.5669999999999999
$ws.Cells.Item(9, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(9, 7).Value = "Benar"

# Row 10: MIB_5.png
$ws.Cells.Item(10, 1).Value = "MIB_5.png"
$ws.Cells.Item(10, 2).Value = 1.278
$ws.Cells.Item(10, 3).Value = 0.004
$ws.Cells.Item(10, 4).Value = 0.973
$ws.Cells.Item(10, 5).Value = 0.633
$ws.Cells.Item(10, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(10, 7).Value = "Benar"

# Row 11: AAH_1.png
$ws.Cells.Item(11, 1).Value = "AAH_1.png"
$ws.Cells.Item(11, 2).Value = 0.852
$ws.Cells.Item(11, 3).Value = 0.003
$ws.Cells.Item(11, 4).Value = 0.97
$ws.Cells.Item(11, 5).Value = 0.433
$ws.Cells.Item(11, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(11, 7).Value = "Salah"

# Row 12: AAH_2.png
$ws.Cells.Item(12, 1).Value = "AAH_2.png"
$ws.Cells.Item(12, 2).Value = 1.141
$ws.Cells.Item(12, 3).Value = 0.004
$ws.Cells.Item(12, 4).Value = 0.989
$ws.Cells.Item(12, 5).Value = 0.7
$ws.Cells.Item(12, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(12, 7).Value = "Benar"

# Row 13: AAH_3.png
$ws.Cells.Item(13, 1).Value = "AAH_3.png"
$ws.Cells.Item(13, 2).Value = 0.851
$ws.Cells.Item(13, 3).Value = 0.003
$ws.Cells.Item(13, 4).Value = 0.967
$ws.Cells.Item(13, 5).Value = 0.633
$ws.Cells.Item(13, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(13, 7).Value = "Benar"

# Row 14: TI_1.png
$ws.Cells.Item(14, 1).Value = "TI_1.png"
$ws.Cells.Item(14, 2).Value = 0.979
$ws.Cells.Item(14, 3).Value = 0.003
$ws.Cells.Item(14, 4).Value = 0.9419999999999999
$ws.Cells.Item(14, 5).Value = 0.367
$ws.Cells.Item(14, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(14, 7).Value = "Salah"

# Row 15: TI_2.png
$ws.Cells.Item(15, 1).Value = "TI_2.png"
$ws.Cells.Item(15, 2).Value = 0.931
$ws.Cells.Item(15, 3).Value = 0.003
$ws.Cells.Item(15, 4).Value = 0.985
$ws.Cells.Item(15, 5).Value = 0.267
$ws.Cells.Item(15, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(15, 7).Value = "Salah"

# Row 16: TI_3.png
$ws.Cells.Item(16, 1).Value = "TI_3.png"
$ws.Cells.Item(16, 2).Value = 0.751
$ws.Cells.Item(16, 3).Value = 0.003
$ws.Cells.Item(16, 4).Value = 0.985
$ws.Cells.Item(16, 5).Value = 0.333
$ws.Cells.Item(16, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(16, 7).Value = "Salah"

# Row 17: TI_4.png
$ws.Cells.Item(17, 1).Value = "TI_4.png"
$ws.Cells.Item(17, 2).Value = 0.74
$ws.Cells.Item(17, 3).Value = 0.002
$ws.Cells.Item(17, 4).Value = 0.958
$ws.Cells.Item(17, 5).Value = 0.467
$ws.Cells.Item(17, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(17, 7).Value = "Salah"

# Row 18: TI_5.png
$ws.Cells.Item(18, 1).Value = "TI_5.png"
$ws.Cells.Item(18, 2).Value = 1.023
$ws.Cells.Item(18, 3).Value = 0.003
$ws.Cells.Item(18, 4).Value = 0.983
$ws.Cells.Item(18, 5).Value = 0.5
$ws.Cells.Item(18, 6).Value = "Toni Ismail"
$ws.Cells.Item(18, 7).Value = "Benar"

# Row 19: RAS_1.png
$ws.Cells.Item(19, 1).Value = "RAS_1.png"
$ws.Cells.Item(19, 2).Value = 0.68
$ws.Cells.Item(19, 3).Value = 0.002
$ws.Cells.Item(19, 4).Value = 0.97
$ws.Cells.Item(19, 5).Value = 0.2
$ws.Cells.Item(19, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(19, 7).Value = "Salah"

# Row 20: RAS_2.png
$ws.Cells.Item(20, 1).Value = "RAS_2.png"
$ws.Cells.Item(20, 2).Value = 1.016
$ws.Cells.Item(20, 3).Value = 0.003
$ws.Cells.Item(20, 4).Value = 0.961
$ws.Cells.Item(20, 5).Value = 0.4
$ws.Cells.Item(20, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(20, 7).Value = "Salah"

# Row 21: RAS_3.png
$ws.Cells.Item(21, 1).Value = "RAS_3.png"
$ws.Cells.Item(21, 2).Value = 0.569
$ws.Cells.Item(21, 3).Value = 0.002
$ws.Cells.Item(21, 4).Value = 0.978
$ws.Cells.Item(21, 5).Value = 0.233
$ws.Cells.Item(21, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(21, 7).Value = "Salah"

# Row 22: RAS_4.png
$ws.Cells.Item(22, 1).Value = "RAS_4.png"
$ws.Cells.Item(22, 2).Value = 1.318
$ws.Cells.Item(22, 3).Value = 0.004
$ws.Cells.Item(22, 4).Value = 0.957
$ws.Cells.Item(22, 5).Value = 0.233
$ws.Cells.Item(22, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(22, 7).Value = "Salah"

# Row 23: RAS_5.png
$ws.Cells.Item(23, 1).Value = "RAS_5.png"
$ws.Cells.Item(23, 2).Value = 1.097
$ws.Cells.Item(23, 3).Value = 0.004
$ws.Cells.Item(23, 4).Value = 0.973
$ws.Cells.Item(23, 5).Value = 0.467
$ws.Cells.Item(23, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(23, 7).Value = "Salah"

# Row 24: RR_1.png
$ws.Cells.Item(24, 1).Value = "RR_1.png"
$ws.Cells.Item(24, 2).Value = 1.155
$ws.Cells.Item(24, 3).Value = 0.004
$ws.Cells.Item(24, 4).Value = 0.985
$ws.Cells.Item(24, 5).Value = 0.367
$ws.Cells.Item(24, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(24, 7).Value = "Salah"

# Row 25: RR_2.png
$ws.Cells.Item(25, 1).Value = "RR_2.png"
$ws.Cells.Item(25, 2).Value = 1.204
$ws.Cells.Item(25, 3).Value = 0.004
$ws.Cells.Item(25, 4).Value = 0.985
$ws.Cells.Item(25, 5).Value = 0.4
$ws.Cells.Item(25, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(25, 7).Value = "Salah"

# Row 26: RR_3.png
$ws.Cells.Item(26, 1).Value = "RR_3.png"
$ws.Cells.Item(26, 2).Value = 0.984
$ws.Cells.Item(26, 3).Value = 0.003
$ws.Cells.Item(26, 4).Value = 0.884
$ws.Cells.Item(26, 5).Value = 0.333
$ws.Cells.Item(26, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(26, 7).Value = "Salah"

# Row 27: RR_4.png
$ws.Cells.Item(27, 1).Value = "RR_4.png"
$ws.Cells.Item(27, 2).Value = 1.185
$ws.Cells.Item(27, 3).Value = 0.004
$ws.Cells.Item(27, 4).Value = 0.99
$ws.Cells.Item(27, 5).Value = 0.467
$ws.Cells.Item(27, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(27, 7).Value = "Salah"

# Row 28: RR_5.png
$ws.Cells.Item(28, 1).Value = "RR_5.png"
$ws.Cells.Item(28, 2).Value = 1.209
$ws.Cells.Item(28, 3).Value = 0.004
$ws.Cells.Item(28, 4).Value = 0.985
$ws.Cells.Item(28, 5).Value = 0.467
$ws.Cells.Item(28, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(28, 7).Value = "Salah"

# Row 29: AR_1.png
$ws.Cells.Item(29, 1).Value = "AR_1.png"
$ws.Cells.Item(29, 2).Value = 0.887
$ws.Cells.Item(29, 3).Value = 0.003
$ws.Cells.Item(29, 4).Value = 0.955
$ws.Cells.Item(29, 5).Value = 0.5
$ws.Cells.Item(29, 6).Value = "Arizli Romadhon"
$ws.Cells.Item(29, 7).Value = "Benar"

# Row 30: GA_1.png
$ws.Cells.Item(30, 1).Value = "GA_1.png"
$ws.Cells.Item(30, 2).Value = 1.209
$ws.Cells.Item(30, 3).Value = 0.004
$ws.Cells.Item(30, 4).Value = 0.971
$ws.Cells.Item(30, 5).Value = 0.733
$ws.Cells.Item(30, 6).Value = "Gege Ardiyansyah"
$ws.Cells.Item(30, 7).Value = "Benar"

# Row 31: GA_2.png
$ws.Cells.Item(31, 1).Value = "GA_2.png"
$ws.Cells.Item(31, 2).Value = 0.6919999999999999
$ws.Cells.Item(31, 3).Value = 0.002
$ws.Cells.Item(31, 4).Value = 0.981
$ws.Cells.Item(31, 5).Value = 0.267
$ws.Cells.Item(31, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(31, 7).Value = "Salah"

# Row 32: GA_3.png
$ws.Cells.Item(32, 1).Value = "GA_3.png"
$ws.Cells.Item(32, 2).Value = 0.773
$ws.Cells.Item(32, 3).Value = 0.003
$ws.Cells.Item(32, 4).Value = 0.981
$ws.Cells.Item(32, 5).Value = 0.3
$ws.Cells.Item(32, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(32, 7).Value = "Salah"

# Row 33: FY_1.png
$ws.Cells.Item(33, 1).Value = "FY_1.png"
$ws.Cells.Item(33, 2).Value = 1.042
$ws.Cells.Item(33, 3).Value = 0.003
$ws.Cells.Item(33, 4).Value = 0.97
$ws.Cells.Item(33, 5).Value = 0.233
$ws.Cells.Item(33, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(33, 7).Value = "Salah"

# Row 34: FY_2.png
$ws.Cells.Item(34, 1).Value = "FY_2.png"
$ws.Cells.Item(34, 2).Value = 1.286
$ws.Cells.Item(34, 3).Value = 0.004
$ws.Cells.Item(34, 4).Value = 0.951
$ws.Cells.Item(34, 5).Value = 0.233
$ws.Cells.Item(34, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(34, 7).Value = "Salah"

# Row 35: FY_3.png
$ws.Cells.Item(35, 1).Value = "FY_3.png"
$ws.Cells.Item(35, 2).Value = 1.246
$ws.Cells.Item(35, 3).Value = 0.004
$ws.Cells.Item(35, 4).Value = 0.985
$ws.Cells.Item(35, 5).Value = 0.233
$ws.Cells.Item(35, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(35, 7).Value = "Salah"

# Row 36: FY_4.png
$ws.Cells.Item(36, 1).Value = "FY_4.png"
$ws.Cells.Item(36, 2).Value = 1.183
$ws.Cells.Item(36, 3).Value = 0.004
$ws.Cells.Item(36, 4).Value = 0.977
$ws.Cells.Item(36, 5).Value = 0.2
$ws.Cells.Item(36, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(36, 7).Value = "Salah"

# Row 37: TO_1.png
$ws.Cells.Item(37, 1).Value = "TO_1.png"
$ws.Cells.Item(37, 2).Value = 0.8159999999999999
$ws.Cells.Item(37, 3).Value = 0.003
$ws.Cells.Item(37, 4).Value = 0.962
$ws.Cells.Item(37, 5).Value = 0.333
$ws.Cells.Item(37, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(37, 7).Value = "Salah"

# Row 38: TO_2.png
$ws.Cells.Item(38, 1).Value = "TO_2.png"
$ws.Cells.Item(38, 2).Value = 0.9360000000000001
$ws.Cells.Item(38, 3).Value = 0.003
$ws.Cells.Item(38, 4).Value = 0.982
$ws.Cells.Item(38, 5).Value = 0.4
$ws.Cells.Item(38, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(38, 7).Value = "Salah"

# Row 39: TO_3.png
$ws.Cells.Item(39, 1).Value = "TO_3.png"
$ws.Cells.Item(39, 2).Value = 0.838
$ws.Cells.Item(39, 3).Value = 0.003
$ws.Cells.Item(39, 4).Value = 0.977
$ws.Cells.Item(39, 5).Value = 0.4
$ws.Cells.Item(39, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(39, 7).Value = "Salah"

# Row 40: TO_4.png
$ws.Cells.Item(40, 1).Value = "TO_4.png"
$ws.Cells.Item(40, 2).Value = 2.492
$ws.Cells.Item(40, 3).Value = 0.008
$ws.Cells.Item(40, 4).Value = 0.783
$ws.Cells.Item(40, 5).Value = 0.333
$ws.Cells.Item(40, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(40, 7).Value = "Salah"

# Row 41: TO_5.png
$ws.Cells.Item(41, 1).Value = "TO_5.png"
$ws.Cells.Item(41, 2).Value = 2.151
$ws.Cells.Item(41, 3).Value = 0.007
$ws.Cells.Item(41, 4).Value = 0.761
$ws.Cells.Item(41, 5).Value = 0.467
$ws.Cells.Item(41, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(41, 7).Value = "Salah"

# Row 42: TD_1.png
$ws.Cells.Item(42, 1).Value = "TD_1.png"
$ws.Cells.Item(42, 2).Value = 1.527
$ws.Cells.Item(42, 3).Value = 0.005
$ws.Cells.Item(42, 4).Value = 0.8100000000000001
$ws.Cells.Item(42, 5).Value = 0.367
$ws.Cells.Item(42, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(42, 7).Value = "Benar"

# Row 43: TD_2.png
$ws.Cells.Item(43, 1).Value = "TD_2.png"
$ws.Cells.Item(43, 2).Value = 1.568
$ws.Cells.Item(43, 3).Value = 0.005
$ws.Cells.Item(43, 4).Value = 0.75
$ws.Cells.Item(43, 5).Value = 0.4
$ws.Cells.Item(43, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(43, 7).Value = "Benar"

# Row 44: TD_3.png
$ws.Cells.Item(44, 1).Value = "TD_3.png"
$ws.Cells.Item(44, 2).Value = 1.04
$ws.Cells.Item(44, 3).Value = 0.003
$ws.Cells.Item(44, 4).Value = 0.944
$ws.Cells.Item(44, 5).Value = 0.2
$ws.Cells.Item(44, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(44, 7).Value = "Benar"

# Row 45: TD_4.png
$ws.Cells.Item(45, 1).Value = "TD_4.png"
$ws.Cells.Item(45, 2).Value = 0.977
$ws.Cells.Item(45, 3).Value = 0.003
$ws.Cells.Item(45, 4).Value = 0.973
$ws.Cells.Item(45, 5).Value = 0.233
$ws.Cells.Item(45, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(45, 7).Value = "Benar"
